$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

# New header row
$ws.Range("A1").Value = "EXPENSE"
$ws.Range("B1").Value = "VALUE"

# The old "Total" label (now on row 6) should be upper-cased.
$ws.Range("A6").Value = "TOTAL"

# Fix the SUM formula to cover the shifted data range (rows 2-5).
$ws.Range("B6").Formula = "=SUM(B2:B5)"
